# #5: property boat&car done
#
# The "汽車" (car) sheet (sheet3) previously only carried the raw
# land/building-style leftover header in row 1 and columns A:G. This
# normalizes row 1 into the same name/capacity/owner/... column header used
# by the other property sheets, and appends the property_category / category
# / date / legislator_name / legislator_id / source_file / index columns
# (H:N) to the two data rows, matching the other sheets (e.g. 土地, 股票).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)   # 汽車

# ---- Row 1: header labels (bold / centered / top-aligned / bordered, like
# the existing B1:G1 header cells) ----
$headers = @{
    2  = "name"
    3  = "capacity"
    4  = "owner"
    5  = "register_date"
    6  = "register_reason"
    7  = "acquire_value"
    8  = "property_category"
    9  = "category"
    10 = "date"
    11 = "legislator_name"
    12 = "legislator_id"
    13 = "source_file"
    14 = "index"
}

foreach ($col in $headers.Keys) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
    $cell.Value2 = $headers[$col]
}

# ---- Rows 2 & 3: fill in the new property_category..index columns (H:N) ----
$row2 = @{
    8  = "land"
    9  = "normal"
    10 = "2012-02-01"
    11 = "吳宜臻"
    12 = 1735
    13 = "tmp2691"
    14 = 32
}
foreach ($col in $row2.Keys) {
    $ws.Cells.Item(2, $col).Value2 = $row2[$col]
}

$row3 = @{
    8  = "land"
    9  = "normal"
    10 = "2012-02-01"
    11 = "吳宜臻"
    12 = 1735
    13 = "tmp2691"
    14 = 33
}
foreach ($col in $row3.Keys) {
    $ws.Cells.Item(3, $col).Value2 = $row3[$col]
}

Write-Host "汽車 sheet header + boat/car metadata columns updated"
